$d = $word.ActiveDocument

# Locate the paragraph we are editing ("When returnbook is called, the
# program should display that patron has incurred fines.") and work out
# absolute character offsets (relative to $d.Content) for the points
# where the run layout needs to change:
#   30 -> right after "...the"      (start of " program...")
#   34 -> right after "...the pro"  (between "pro" and "gram")
#   84 -> right after "...fines"    (start of the trailing ".")
$text = $d.Content.Text
$paraStart = $text.IndexOf("When returnbook")

$boundaryTheEnd = $paraStart + 30
$boundaryProEnd = $paraStart + 34
$boundaryFinesEnd = $paraStart + 84

# Drop temporary bookmarks at those points first. A bookmark splits the
# run it lands in, and - so long as it stays in place - also stops the
# Find/Replace engine from coalescing runs across it, which lets each
# Find/Replace below only touch the text it is meant to.
$d.Bookmarks.Add("ZZTmpBoundary1", $d.Range($boundaryTheEnd, $boundaryTheEnd))
$d.Bookmarks.Add("ZZTmpBoundary2", $d.Range($boundaryProEnd, $boundaryProEnd))
$d.Bookmarks.Add("ZZTmpBoundary3", $d.Range($boundaryFinesEnd, $boundaryFinesEnd))

# --------------------------------------------------------------------
# 1) "When " + "returnbook" (spell-checked) + " is called, the" become
#    a single plain run "When returnbook is called, the" - the
#    proofErr spell-check markers disappear as a side effect of the
#    text edit. Stops at ZZTmpBoundary1, leaving " program..." alone.
# --------------------------------------------------------------------
$d.Content.Find.Execute("When returnbook is called, the", $true, $false, $false, $false, $false, $true, 1, $false, "When returnbook is called, the", 2)

# --------------------------------------------------------------------
# 2) Touch " pro" (no textual change) so it becomes its own clean run
#    between ZZTmpBoundary1 and ZZTmpBoundary2.
# --------------------------------------------------------------------
$d.Content.Find.Execute(" pro", $true, $false, $false, $false, $false, $true, 1, $false, " pro", 2)

# --------------------------------------------------------------------
# 3) Type a stray "z" right before "gram...": "program" becomes
#    "proz" + "gram". Bounded by ZZTmpBoundary2 .. ZZTmpBoundary3, so
#    it neither pulls in " pro" nor the trailing ".".
# --------------------------------------------------------------------
$d.Content.Find.Execute("gram should display that patron has incurred fines", $true, $false, $false, $false, $false, $true, 1, $false, "zgram should display that patron has incurred fines", 2)

# --------------------------------------------------------------------
# 4) Drop the real "_GoBack" bookmark exactly where the "z" was typed
#    (between "z" and "gram"), splitting that run into "z" / "gram...".
#    Word only ever keeps a single "_GoBack" bookmark tracking the most
#    recent edit location, so adding it here automatically removes the
#    old one (in the "Setup" bullet below) and renumbers whichever
#    bookmark had claimed its freed-up id.
# --------------------------------------------------------------------
$goBackPos = $paraStart + 35
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos))

# --------------------------------------------------------------------
# 5) Remove the temporary boundary bookmarks - the run splits they
#    created persist even once the bookmarks themselves are gone.
# --------------------------------------------------------------------
$d.Bookmarks("ZZTmpBoundary1").Delete()
$d.Bookmarks("ZZTmpBoundary2").Delete()
$d.Bookmarks("ZZTmpBoundary3").Delete()
